# "cleaned up all scripts" - rename the decimal lat/lon header labels
# from lat_dec/lon_dec to lat/lon (columns D and E on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D1").Value = "lat"
$ws.Range("E1").Value = "lon"

# Leave the selection on E2 rather than the stale E12.
$ws.Range("E2").Select()

$wb.Save()
